$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = 67
$ws.Range("B68").Value = 1
$ws.Range("C68").Value = "2024-06-16 08:14:36"
$ws.Range("D68").Value = 200
$ws.Range("E68").Value = 3

$ws.Range("A69").Value = 68
$ws.Range("B69").Value = 2
$ws.Range("C69").Value = "2024-06-16 08:14:37"
$ws.Range("D69").Value = 200
$ws.Range("E69").Value = 0
